$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'315.85"
$ws.Range("E2").Value = "'3.48%"
$ws.Range("D3").Value = "'35.60"
$ws.Range("E3").Value = "'-0.90%"
$ws.Range("D4").Value = "'5.123"
$ws.Range("E4").Value = "'0.81%"
$ws.Range("D5").Value = "'0.08104"
$ws.Range("E5").Value = "'3.02%"
$ws.Range("E6").Value = "'0.61%"
$ws.Range("D7").Value = "'8.014"
$ws.Range("E7").Value = "'1.28%"
$ws.Range("D8").Value = "'4.146"
$ws.Range("E8").Value = "'1.11%"
$ws.Range("D9").Value = "'0.9269"
$ws.Range("E9").Value = "'0.89%"
$ws.Range("D10").Value = "'0.1018"
$ws.Range("E10").Value = "'5.15%"
$ws.Range("D11").Value = "'0.1871"
$ws.Range("E11").Value = "'0.93%"
$ws.Range("D12").Value = "'0.09213"
$ws.Range("E12").Value = "'7.07%"
$ws.Range("D13").Value = "'0.03601"
$ws.Range("E13").Value = "'0.74%"
$ws.Range("D14").Value = "'0.09889"
$ws.Range("E14").Value = "'-0.45%"
$ws.Range("D15").Value = "'0.001445"
$ws.Range("E15").Value = "'1.43%"
$ws.Range("D16").Value = "'0.005760"
$ws.Range("E16").Value = "'1.79%"
$ws.Range("D17").Value = "'3.460"
$ws.Range("E17").Value = "'0.63%"
$ws.Range("D18").Value = "'2.735"
$ws.Range("E18").Value = "'3.26%"
$ws.Range("D19").Value = "'0.3366"
$ws.Range("E19").Value = "'-0.94%"
$ws.Range("D20").Value = "'0.1333"
$ws.Range("E20").Value = "'1.23%"
$ws.Range("D21").Value = "'5.137"
$ws.Range("E21").Value = "'-0.21%"
$ws.Range("D22").Value = "'0.2223"
$ws.Range("E22").Value = "'1.03%"
$ws.Range("D23").Value = "'0.04582"
$ws.Range("E23").Value = "'0.92%"
$ws.Range("E24").Value = "'1.21%"
$ws.Range("D25").Value = "'0.004706"
$ws.Range("E25").Value = "'-6.97%"
$ws.Range("D26").Value = "'0.0001251"
$ws.Range("E26").Value = "'-21.94%"
$ws.Range("D27").Value = "'0.0004504"
$ws.Range("E27").Value = "'-4.91%"
$ws.Range("D39").Value = "'0.01949"
$ws.Range("E39").Value = "'5.64%"
$ws.Range("D40").Value = "'0.04862"
$ws.Range("E40").Value = "'2.54%"
$ws.Range("D41").Value = "'0.007759"
$ws.Range("E41").Value = "'2.97%"
$ws.Range("D42").Value = "'0.1388"
$ws.Range("E42").Value = "'-0.79%"
$ws.Range("D43").Value = "'0.007699"
$ws.Range("E43").Value = "'-0.06%"
$ws.Range("D44").Value = "'0.002105"
$ws.Range("E44").Value = "'-4.55%"
$ws.Range("D45").Value = "'0.01162"
$ws.Range("E45").Value = "'5.49%"
$ws.Range("D46").Value = "'0.00006495"
$ws.Range("E46").Value = "'2.71%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.36%"
$ws.Range("D48").Value = "'39.18"
$ws.Range("E48").Value = "'-17.45%"
$ws.Range("D49").Value = "'0.001702"
$ws.Range("E49").Value = "'-14.69%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.36%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.36%"
